# Widen the four step text boxes on slide 15 ("Ipv4_match", "Check_ttl",
# "Set_dmac", "Set_smac") to make room for the new packet-filter figure.
#
# NOTE: PowerPoint's Shape.Left/.Top/.Width/.Height properties are stored
# internally as single-precision (f32) point values, and converting a point
# value to EMU truncates (floors) rather than rounds. The literal point
# values below were chosen so that, after that float32 truncation, they
# reproduce the exact target EMU offsets/extents.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# "Ipv4_match" text box (off x 1080229 -> 1080228, cx 1339250 -> 1509535)
$sh = $s.Shapes.Item(6)
$sh.Width = 118.86102682204724
$sh.Left = 85.05732353464568

# "Check_ttl" text box (off x unchanged 3019156, cx 1237992 -> 1447246)
$sh = $s.Shapes.Item(7)
$sh.Width = 113.9563789527559
$sh.Left = 237.7288131976378

# "Set_dmac" text box (off x 4789292 -> 4789291, cx 1177812 -> 1369383)
$sh = $s.Shapes.Item(8)
$sh.Width = 107.82543567086614
$sh.Left = 377.1095123590551

# "Set_smac" text box (off x 6634055 -> 6563102, cx 1113955 -> 1362290)
$sh = $s.Shapes.Item(9)
$sh.Width = 107.26692583385827
$sh.Left = 516.7796936393701
